# Apple Numbers export compatibility: add two new rows to the "Empty rows
# and columns" sheet, widen the columns that now hold the longer strings,
# and make that sheet the active tab (it was "Basic types" before).

$wb = $excel.ActiveWorkbook

$wsEmpty = $wb.Worksheets.Item("Empty rows and columns")

# New data rows (row 3 stays blank, as before; rows 5 and 6 are new).
$wsEmpty.Range("B5").Value = "Row with leading null"
$wsEmpty.Range("D5").Value = "Row with leading null"

$wsEmpty.Range("A6").Value = "Row with trailing null"
$wsEmpty.Range("B6").Value = "Row with trailing null"

# Widen the columns that now contain the longer strings above.
$wsEmpty.Columns.Item(1).ColumnWidth = 19.333333333333332
$wsEmpty.Columns.Item(2).ColumnWidth = 19.666666666666668
$wsEmpty.Columns.Item(4).ColumnWidth = 19.666666666666668

# "Empty rows and columns" becomes the active/selected sheet & tab
# (previously "Basic types" was tabSelected / the active tab).
$wsEmpty.Activate()
